$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76, pushing existing rows 76..187 down to 77..188
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new record
$ws.Range("A76").Value = 5
$ws.Range("B76").Value = "Macroferia Regional de Talca"
$ws.Range("C76").Value = "Maule"
$ws.Range("D76").Value = 44557
$ws.Range("E76").Value = 7
$ws.Range("F76").Value = 100112008
$ws.Range("G76").Value = "Coliflor"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 900
$ws.Range("L76").Value = 900
$ws.Range("M76").Value = 900
$ws.Range("N76").Value = "`$/unidad"
$ws.Range("O76").Value = "Región del Maule"
$ws.Range("P76").Value = 900
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"
